$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.648.50"
$ws.Range("E2").Value = "  +2.80%  "

# Row 3
$ws.Range("D3").Value = "3.202.10"
$ws.Range("E3").Value = "  +1.51%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").Value = "598.79"
$ws.Range("E5").Value = "  +3.30%  "

# Row 6
$ws.Range("D6").Value = "156.06"
$ws.Range("E6").Value = "  +4.41%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "0.560"
$ws.Range("E8").Value = "  +6.27%  "

# Row 9
$ws.Range("D9").Value = "3.200.11"
$ws.Range("E9").Value = "  +1.28%  "

# Row 10
$ws.Range("E10").Value = "  +1.34%  "

# Row 11
$ws.Range("E11").Value = "  -4.55%  "

# Row 12
$ws.Range("D12").Value = "0.522"
$ws.Range("E12").Value = "  +3.73%  "

# Row 13
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +2.02%  "

# Row 14
$ws.Range("D14").Value = "39.41"
$ws.Range("E14").Value = "  +5.65%  "

# Row 15
$ws.Range("D15").Value = "3.724.59"
$ws.Range("E15").Value = "  +1.81%  "

# Row 16
$ws.Range("E16").Value = "  +4.79%  "

# Row 17
$ws.Range("D17").Value = "66.641.29"
$ws.Range("E17").Value = "  +2.82%  "

# Row 18
$ws.Range("D18").Value = "3.203.71"
$ws.Range("E18").Value = "  +1.60%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "525.08"
$ws.Range("E19").Value = "  +3.55%  "

# Row 20
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "0.112"
$ws.Range("E20").Value = "  +0.73%  "

# Row 21
$ws.Range("D21").Value = "15.50"
$ws.Range("E21").Value = "  +3.80%  "

# Row 22
$ws.Range("E22").Value = "  +3.44%  "

# Row 23
$ws.Range("D23").Value = "8.18"
$ws.Range("E23").Value = "  +5.51%  "

# Row 24
$ws.Range("D24").Value = "15.04"
$ws.Range("E24").Value = "  -1.12%  "

# Row 25
$ws.Range("D25").Value = "85.98"
$ws.Range("E25").Value = "  +1.62%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").Value = "9.30"
$ws.Range("E27").Value = "  +3.26%  "

# Row 28
$ws.Range("E28").Value = "  +3.15%  "

# Row 29
$ws.Range("D29").Value = "2.38"
$ws.Range("E29").Value = "  +8.52%  "

# Row 30
$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").Value = "2.97"
$ws.Range("E30").Value = "  +6.69%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "7.06"
$ws.Range("E31").Value = "  +12.07%  "

# Row 32
$ws.Range("D32").Value = "28.41"
$ws.Range("E32").Value = "  +2.50%  "

# Row 33
$ws.Range("D33").Value = "1.24"
$ws.Range("E33").Value = "  +3.10%  "

# Row 34
$ws.Range("E34").Value = "  +0.23%  "

# Row 35
$ws.Range("D35").Value = "6.58"
$ws.Range("E35").Value = "  +1.09%  "

# Row 36
$ws.Range("D36").Value = "510.38"
$ws.Range("E36").Value = "  +6.38%  "

# Row 37
$ws.Range("D37").Value = "54.95"
$ws.Range("E37").Value = "  +0.11%  "

# Row 38
$ws.Range("D38").Value = "0.0909"
$ws.Range("E38").Value = "  +1.65%  "

# Row 39
$ws.Range("D39").Value = "0.0427"
$ws.Range("E39").Value = "  +2.77%  "

# Row 40
$ws.Range("E40").Value = "  +9.39%  "

# Row 41
$ws.Range("E41").Value = "  +1.92%  "

# Row 42
$ws.Range("E42").Value = "  -0.76%  "

# Row 43
$ws.Range("E43").Value = "  +16.08%  "

# Row 44
$ws.Range("E44").Value = "  +6.97%  "

# Row 45
$ws.Range("E45").Value = "  +1.14%  "

# Row 46
$ws.Range("D46").Value = "2.902.19"
$ws.Range("E46").Value = "  -3.07%  "

# Row 47
$ws.Range("D47").Value = "28.60"
$ws.Range("E47").Value = "  +0.93%  "

# Row 48
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  +10.68%  "

# Row 49
$ws.Range("E49").Value = "  +3.57%  "

# Row 50
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "2.37"
$ws.Range("E50").Value = "  +5.43%  "

# Row 51
$ws.Range("B51").Value = "USDe"
$ws.Range("C51").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.04%  "
